$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The account-statement detail block (rows 16-36) is reorganised so the data
# is grouped by period (ascending 1607..1701) and, within each period, lists
# the three workers (Natalia, Lizeth, Mayury) in order - this also folds in
# the newly added Lizeth/Mayury periods that used to be appended after Natalia.
$data = @(
    ,@(16, "64586905", "NATALIA INES RODRIGUEZ JEREZ", "1607", 76000, 1900000)
    ,@(17, "1052972861", "LIZETH PAOLA ARANGO VARGAS", "1607", 65880, 1647000)
    ,@(18, "455607942", "MAYURY MARIA SIERRA GREY", "1607", 24640, 781242)
    ,@(19, "64586905", "NATALIA INES RODRIGUEZ JEREZ", "1608", 76000, 1900000)
    ,@(20, "1052972861", "LIZETH PAOLA ARANGO VARGAS", "1608", 65880, 1647000)
    ,@(21, "455607942", "MAYURY MARIA SIERRA GREY", "1608", 24640, 781242)
    ,@(22, "64586905", "NATALIA INES RODRIGUEZ JEREZ", "1609", 76000, 1900000)
    ,@(23, "1052972861", "LIZETH PAOLA ARANGO VARGAS", "1609", 65880, 1647000)
    ,@(24, "455607942", "MAYURY MARIA SIERRA GREY", "1609", 24640, 781242)
    ,@(25, "64586905", "NATALIA INES RODRIGUEZ JEREZ", "1610", 76000, 1900000)
    ,@(26, "1052972861", "LIZETH PAOLA ARANGO VARGAS", "1610", 65880, 1647000)
    ,@(27, "455607942", "MAYURY MARIA SIERRA GREY", "1610", 24640, 781242)
    ,@(28, "64586905", "NATALIA INES RODRIGUEZ JEREZ", "1611", 76000, 1900000)
    ,@(29, "1052972861", "LIZETH PAOLA ARANGO VARGAS", "1611", 65880, 1647000)
    ,@(30, "455607942", "MAYURY MARIA SIERRA GREY", "1611", 24640, 781242)
    ,@(31, "64586905", "NATALIA INES RODRIGUEZ JEREZ", "1612", 76000, 1900000)
    ,@(32, "1052972861", "LIZETH PAOLA ARANGO VARGAS", "1612", 65880, 1647000)
    ,@(33, "455607942", "MAYURY MARIA SIERRA GREY", "1612", 24640, 781242)
    ,@(34, "64586905", "NATALIA INES RODRIGUEZ JEREZ", "1701", 76000, 1900000)
    ,@(35, "1052972861", "LIZETH PAOLA ARANGO VARGAS", "1701", 65880, 1647000)
    ,@(36, "455607942", "MAYURY MARIA SIERRA GREY", "1701", 24640, 781242)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
